$wb = $excel.ActiveWorkbook

# Sheets "展览" (exhibitions) and "全部类型" (all types) carry duplicate data
# tables; update the "想去人数" (people-who-want-to-go) counts in both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1426
    $ws.Range("F5").Value = 10
    $ws.Range("F8").Value = 198
}
